# The sheet originally held a 2-column "productLine/totalSales" table in
# A1:B8. The new data is a 5-column "title/tactic/technique/other_info/
# description" table in A1:E6, so first drop the old rows (7:8) that fall
# outside the new range before writing the replacement content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:B8").ClearContents()

# Header row (A1:B1 already carry the bold/bordered header style from the
# original sheet; copy that formatting across to the newly-added C1:E1
# header cells so the whole header row looks consistent).
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "tactic"
$ws.Range("C1").Value = "technique"
$ws.Range("D1").Value = "other_info"
$ws.Range("E1").Value = "description"

$ws.Range("A1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

# Data rows
$data = @(
    @("Antivirus Exploitation Framework Detection", "Execution", "T1203", "", "Exploitation for Client Execution : Detecting exploitation of vulnerabilities for execution."),
    @("Antivirus Hacktool Detection", "Execution", "T1204.002", "", "User Execution: Malicious File : Detection of tools requiring user execution to compromise systems."),
    @("Antivirus Password Dumper Detection", "Credential Access", "T1003", "", "OS Credential Dumping : Detecting tools attempting to dump OS credentials from memory."),
    @("Antivirus PrinterNightmare CVE-2021-34527 Exploit Detection", "Privilege Escalation", "T1068", "", "Exploitation for Privilege Escalation : Detection of PrinterNightmare (CVE-2021-34527) exploit for privilege escalation."),
    @("Antivirus Ransomware Detection", "Impact", "T1486", "", "Data Encrypted for Impact : Detection of ransomware encrypting files to cause impact.")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}
